# Junction_Flooding_51.xlsx edit:
#  - Row 5 values rounded to "custom accuracy" (2 decimal places)
#  - Row 6 (the last data row) removed entirely
#  - sheet dimension shrinks from A1:AH6 to A1:AH5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Overwrite row 5's numeric values with their 2-decimal-place equivalents.
$ws.Range("B5").Value  = 9.95
$ws.Range("C5").Value  = 7.41
$ws.Range("D5").Value  = 0.86
$ws.Range("E5").Value  = 21.92
$ws.Range("F5").Value  = 17.59
$ws.Range("G5").Value  = 7.71
$ws.Range("H5").Value  = 30.81
$ws.Range("I5").Value  = 12.16
$ws.Range("J5").Value  = 5.49
$ws.Range("K5").Value  = 7.67
$ws.Range("L5").Value  = 8.81
$ws.Range("M5").Value  = 9.41
$ws.Range("N5").Value  = 2.54
$ws.Range("O5").Value  = 7.9
$ws.Range("P5").Value  = 11.09
$ws.Range("Q5").Value  = 6.84
$ws.Range("R5").Value  = 0.64
$ws.Range("S5").Value  = 0.39
$ws.Range("T5").Value  = 112.84
$ws.Range("U5").Value  = 22.07
$ws.Range("V5").Value  = 7.29
$ws.Range("W5").Value  = 14.65
$ws.Range("X5").Value  = 7.72
$ws.Range("Y5").Value  = 1.06
$ws.Range("Z5").Value  = 15.36
$ws.Range("AA5").Value = 6.44
$ws.Range("AB5").Value = 5.78
$ws.Range("AC5").Value = 6.79
$ws.Range("AD5").Value = 9.28
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 28.1
$ws.Range("AG5").Value = 4
$ws.Range("AH5").Value = 9.11

# 2) Delete the now-redundant last data row (old row 6), shifting nothing
#    up underneath it since it was the final row; this also shrinks the
#    sheet's used range / dimension down to A1:AH5.
$ws.Rows("6").Delete()
